# Update per-row profit/price figures on each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed market-board pricing pulled in by the scheduled Sheets runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 11.25  # H11
$ws.Cells.Item(11, 9).Value = 11.25  # I11
$ws.Cells.Item(11, 11).Value = 11.25  # K11
$ws.Cells.Item(11, 13).Value = 128.75  # M11
$ws.Cells.Item(15, 8).Value = 730.5625  # H15
$ws.Cells.Item(15, 9).Value = 730.5625  # I15
$ws.Cells.Item(15, 11).Value = 2191.6875  # K15
$ws.Cells.Item(15, 13).Value = -2022.6875  # M15
$ws.Cells.Item(28, 8).Value = 221.28572  # H28
$ws.Cells.Item(28, 9).Value = 91.5  # I28
$ws.Cells.Item(28, 10).Value = 1000  # J28
$ws.Cells.Item(28, 11).Value = 91.5  # K28
$ws.Cells.Item(28, 12).Value = 1000  # L28
$ws.Cells.Item(28, 13).Value = 393.5  # M28
$ws.Cells.Item(28, 14).Value = -1970  # N28
$ws.Cells.Item(33, 8).Value = 204.55  # H33
$ws.Cells.Item(33, 9).Value = 204.78947  # I33
$ws.Cells.Item(33, 11).Value = 204.78947  # K33
$ws.Cells.Item(33, 13).Value = 24.21053000000001  # M33
$ws.Cells.Item(92, 8).Value = 252.8  # H92
$ws.Cells.Item(92, 9).Value = 252.8  # I92
$ws.Cells.Item(92, 11).Value = 252.8  # K92
$ws.Cells.Item(92, 13).Value = 995.2  # M92
$ws.Cells.Item(96, 8).Value = 767.75  # H96
$ws.Cells.Item(96, 9).Value = 87.5  # I96
$ws.Cells.Item(96, 11).Value = 262.5  # K96
$ws.Cells.Item(96, 13).Value = 1110.5  # M96
$ws.Cells.Item(116, 8).Value = 3747.5  # H116
$ws.Cells.Item(116, 10).Value = 3500  # J116
$ws.Cells.Item(116, 12).Value = 3500  # L116
$ws.Cells.Item(116, 14).Value = -10384  # N116
$ws.Cells.Item(125, 8).Value = 500  # H125
$ws.Cells.Item(125, 9).Value = 500  # I125
$ws.Cells.Item(125, 11).Value = 4500  # K125
$ws.Cells.Item(125, 13).Value = -2040  # M125
$ws.Cells.Item(132, 8).Value = 3184.5715  # H132
$ws.Cells.Item(132, 9).Value = 3048.6667  # I132
$ws.Cells.Item(132, 10).Value = 4000  # J132
$ws.Cells.Item(132, 11).Value = 9146.000100000001  # K132
$ws.Cells.Item(132, 12).Value = 12000  # L132
$ws.Cells.Item(132, 13).Value = -6616.000100000001  # M132
$ws.Cells.Item(132, 14).Value = -17060  # N132
$ws.Cells.Item(138, 8).Value = 1880  # H138
$ws.Cells.Item(138, 9).Value = 1000  # I138
$ws.Cells.Item(138, 10).Value = 2466.6667  # J138
$ws.Cells.Item(138, 11).Value = 3000  # K138
$ws.Cells.Item(138, 12).Value = 7400.000100000001  # L138
$ws.Cells.Item(138, 13).Value = 2140  # M138
$ws.Cells.Item(138, 14).Value = -17680.0001  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 670  # H12
$ws.Cells.Item(12, 9).Value = 450  # I12
$ws.Cells.Item(12, 10).Value = 1000  # J12
$ws.Cells.Item(12, 11).Value = 450  # K12
$ws.Cells.Item(12, 12).Value = 1000  # L12
$ws.Cells.Item(12, 13).Value = -277  # M12
$ws.Cells.Item(12, 14).Value = -1346  # N12
$ws.Cells.Item(61, 8).Value = 4000  # H61
$ws.Cells.Item(61, 9).Value = 0  # I61
$ws.Cells.Item(61, 10).Value = 4000  # J61
$ws.Cells.Item(61, 11).Value = 0  # K61
$ws.Cells.Item(61, 13).Value = 4000  # M61
$ws.Cells.Item(61, 14).Value = -4424  # N61
$ws.Cells.Item(103, 8).Value = 10787  # H103
$ws.Cells.Item(103, 10).Value = 10787  # J103
$ws.Cells.Item(103, 12).Value = 10787  # L103
$ws.Cells.Item(103, 14).Value = -13131  # N103
$ws.Cells.Item(110, 8).Value = 1007.3333  # H110
$ws.Cells.Item(110, 9).Value = 911  # I110
$ws.Cells.Item(110, 10).Value = 1055.5  # J110
$ws.Cells.Item(110, 11).Value = 911  # K110
$ws.Cells.Item(110, 12).Value = 1055.5  # L110
$ws.Cells.Item(110, 13).Value = 1134  # M110
$ws.Cells.Item(110, 14).Value = -5145.5  # N110
$ws.Cells.Item(136, 8).Value = 4000  # H136
$ws.Cells.Item(136, 9).Value = 0  # I136
$ws.Cells.Item(136, 10).Value = 4000  # J136
$ws.Cells.Item(136, 11).Value = 0  # K136
$ws.Cells.Item(136, 13).Value = 12000  # M136
$ws.Cells.Item(136, 14).Value = -17100  # N136
$ws.Cells.Item(141, 8).Value = 195000  # H141
$ws.Cells.Item(141, 10).Value = 195000  # J141
$ws.Cells.Item(141, 12).Value = 195000  # L141
$ws.Cells.Item(141, 14).Value = -205360  # N141

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1851  # H20
$ws.Cells.Item(20, 9).Value = 1302.875  # I20
$ws.Cells.Item(20, 11).Value = 1302.875  # K20
$ws.Cells.Item(20, 13).Value = -1055.875  # M20
$ws.Cells.Item(80, 8).Value = 938.8333  # H80
$ws.Cells.Item(80, 10).Value = 977  # J80
$ws.Cells.Item(80, 12).Value = 977  # L80
$ws.Cells.Item(80, 14).Value = -2973  # N80
$ws.Cells.Item(83, 8).Value = 938.8333  # H83
$ws.Cells.Item(83, 10).Value = 977  # J83
$ws.Cells.Item(83, 12).Value = 4885  # L83
$ws.Cells.Item(83, 14).Value = -14869  # N83
$ws.Cells.Item(86, 8).Value = 1363  # H86
$ws.Cells.Item(86, 9).Value = 1455.1666  # I86
$ws.Cells.Item(86, 10).Value = 1224.75  # J86
$ws.Cells.Item(86, 11).Value = 1455.1666  # K86
$ws.Cells.Item(86, 12).Value = 1224.75  # L86
$ws.Cells.Item(86, 13).Value = -332.1666  # M86
$ws.Cells.Item(86, 14).Value = -3470.75  # N86
$ws.Cells.Item(89, 8).Value = 1363  # H89
$ws.Cells.Item(89, 9).Value = 1455.1666  # I89
$ws.Cells.Item(89, 10).Value = 1224.75  # J89
$ws.Cells.Item(89, 11).Value = 7275.833000000001  # K89
$ws.Cells.Item(89, 12).Value = 6123.75  # L89
$ws.Cells.Item(89, 13).Value = -1659.833000000001  # M89
$ws.Cells.Item(89, 14).Value = -17355.75  # N89
$ws.Cells.Item(103, 8).Value = 12185  # H103
$ws.Cells.Item(103, 10).Value = 12185  # J103
$ws.Cells.Item(103, 12).Value = 12185  # L103
$ws.Cells.Item(103, 14).Value = -14529  # N103

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(43, 8).Value = 22871.334  # H43
$ws.Cells.Item(43, 10).Value = 22871.334  # J43
$ws.Cells.Item(43, 12).Value = 22871.334  # L43
$ws.Cells.Item(43, 14).Value = -23239.334  # N43
$ws.Cells.Item(101, 8).Value = 22871.334  # H101
$ws.Cells.Item(101, 10).Value = 22871.334  # J101
$ws.Cells.Item(101, 12).Value = 22871.334  # L101
$ws.Cells.Item(101, 14).Value = -29361.334  # N101
$ws.Cells.Item(132, 8).Value = 4833.3335  # H132
$ws.Cells.Item(132, 9).Value = 4833.3335  # I132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 11).Value = 14500.0005  # K132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 14).Value = -11970.0005  # N132
$ws.Cells.Item(141, 8).Value = 0  # H141
$ws.Cells.Item(141, 10).Value = 0  # J141
$ws.Cells.Item(141, 14).Value = 0  # N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 410.375  # H4
$ws.Cells.Item(4, 9).Value = 410.375  # I4
$ws.Cells.Item(4, 11).Value = 1231.125  # K4
$ws.Cells.Item(4, 13).Value = -1119.125  # M4
$ws.Cells.Item(51, 8).Value = 700  # H51
$ws.Cells.Item(51, 10).Value = 625  # J51
$ws.Cells.Item(51, 12).Value = 1875  # L51
$ws.Cells.Item(51, 14).Value = -2795  # N51
$ws.Cells.Item(101, 8).Value = 0  # H101
$ws.Cells.Item(101, 9).Value = 0  # I101
$ws.Cells.Item(101, 11).Value = 0  # K101
$ws.Cells.Item(104, 8).Value = 3725.5  # H104
$ws.Cells.Item(104, 9).Value = 3725.5  # I104
$ws.Cells.Item(104, 10).Value = 0  # J104
$ws.Cells.Item(104, 11).Value = 11176.5  # K104
$ws.Cells.Item(104, 12).Value = 0  # L104
$ws.Cells.Item(104, 14).Value = -8555.5  # N104
$ws.Cells.Item(106, 8).Value = 0  # H106
$ws.Cells.Item(106, 10).Value = 0  # J106
$ws.Cells.Item(106, 14).Value = 0  # N106

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 14998.667  # H24
$ws.Cells.Item(24, 10).Value = 14998.667  # J24
$ws.Cells.Item(24, 12).Value = 14998.667  # L24
$ws.Cells.Item(24, 14).Value = -15344.667  # N24
$ws.Cells.Item(70, 8).Value = 5374.5  # H70
$ws.Cells.Item(70, 9).Value = 4250  # I70
$ws.Cells.Item(70, 11).Value = 4250  # K70
$ws.Cells.Item(70, 13).Value = -3980  # M70
$ws.Cells.Item(73, 8).Value = 5374.5  # H73
$ws.Cells.Item(73, 9).Value = 4250  # I73
$ws.Cells.Item(73, 11).Value = 4250  # K73
$ws.Cells.Item(73, 13).Value = -3314  # M73
$ws.Cells.Item(102, 8).Value = 0  # H102
$ws.Cells.Item(102, 9).Value = 0  # I102
$ws.Cells.Item(102, 11).Value = 0  # K102
$ws.Cells.Item(132, 8).Value = 3112.5  # H132
$ws.Cells.Item(132, 9).Value = 2225  # I132
$ws.Cells.Item(132, 10).Value = 4000  # J132
$ws.Cells.Item(132, 11).Value = 6675  # K132
$ws.Cells.Item(132, 12).Value = 12000  # L132
$ws.Cells.Item(132, 13).Value = -4145  # M132
$ws.Cells.Item(132, 14).Value = -17060  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 0  # H22
$ws.Cells.Item(22, 10).Value = 0  # J22
$ws.Cells.Item(22, 14).Value = 0  # N22
$ws.Cells.Item(27, 8).Value = 0  # H27
$ws.Cells.Item(27, 10).Value = 0  # J27
$ws.Cells.Item(27, 14).Value = 0  # N27
$ws.Cells.Item(132, 8).Value = 4294  # H132
$ws.Cells.Item(132, 10).Value = 2600  # J132
$ws.Cells.Item(132, 12).Value = 7800  # L132
$ws.Cells.Item(132, 14).Value = -12860  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(63, 8).Value = 18462.5  # H63
$ws.Cells.Item(63, 10).Value = 18462.5  # J63
$ws.Cells.Item(63, 12).Value = 18462.5  # L63
$ws.Cells.Item(63, 14).Value = -19710.5  # N63
$ws.Cells.Item(66, 8).Value = 18462.5  # H66
$ws.Cells.Item(66, 10).Value = 18462.5  # J66
$ws.Cells.Item(66, 12).Value = 55387.5  # L66
$ws.Cells.Item(66, 14).Value = -61627.5  # N66
$ws.Cells.Item(69, 8).Value = 27592.2  # H69
$ws.Cells.Item(69, 10).Value = 27592.2  # J69
$ws.Cells.Item(69, 12).Value = 27592.2  # L69
$ws.Cells.Item(69, 14).Value = -29090.2  # N69
$ws.Cells.Item(72, 8).Value = 27592.2  # H72
$ws.Cells.Item(72, 10).Value = 27592.2  # J72
$ws.Cells.Item(72, 12).Value = 82776.60000000001  # L72
$ws.Cells.Item(72, 14).Value = -90264.60000000001  # N72
$ws.Cells.Item(95, 8).Value = 0  # H95
$ws.Cells.Item(95, 10).Value = 0  # J95
$ws.Cells.Item(95, 14).Value = 0  # N95
$ws.Cells.Item(132, 8).Value = 0  # H132
$ws.Cells.Item(132, 9).Value = 0  # I132
$ws.Cells.Item(132, 11).Value = 0  # K132
